# Add updated microstates from SAMPL6
# Adds a new "canonical SMILES" column (D) next to the existing
# "canonical isomeric SMILES" column (C). The canonical SMILES is simply
# the isomeric SMILES with the E/Z stereo-bond markers ("/" and "\")
# stripped out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column D, on the same header row as the other columns.
$ws.Range("D2").Value = "canonical SMILES"

# For every data row, derive the canonical SMILES from the canonical
# isomeric SMILES already present in column C and write it to column D.
for ($r = 3; $r -le 17; $r++) {
    $isomeric = $ws.Cells.Item($r, 3).Value()
    if ($isomeric -ne $null) {
        $canonical = $isomeric.Replace("/", "").Replace("\", "")
        $ws.Cells.Item($r, 4).Value = $canonical
    }
}

# Give the new column a sensible width, matching the other SMILES columns.
$ws.Columns.Item(4).ColumnWidth = 42.28515625
